$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.108069666666667
$ws.Range("H2").Value = 27.324209
$ws.Range("I2").Value = 0.00155006418458712
$ws.Range("J2").Value = 0.00155006418458712
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.22896466666667
$ws.Range("N2").Value = 93.686894
$ws.Range("O2").Value = 0.2877106972998646
$ws.Range("P2").Value = 0.2877106972998646
$ws.Range("Q2").Value = 284.4355858018718
$ws.Range("R2").Value = 2559.920272216846
$ws.Range("S2").Value = 0.0004459700474071062
$ws.Range("T2").Value = 0.0004459700474071062

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.108069666666667
$ws.Range("H3").Value = 27.324209
$ws.Range("I3").Value = 0.00155006418458712
$ws.Range("J3").Value = 0.00155006418458712
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.44578266666667
$ws.Range("N3").Value = 121.337348
$ws.Range("O3").Value = 0.3726247238124506
$ws.Range("P3").Value = 0.3726247238124505
$ws.Range("Q3").Value = 368.3830062508591
$ws.Range("R3").Value = 3315.447056257732
$ws.Range("S3").Value = 0.0005775922386733469
$ws.Range("T3").Value = 0.0005775922386733467

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.108069666666667
$ws.Range("H4").Value = 27.324209
$ws.Range("I4").Value = 0.00155006418458712
$ws.Range("J4").Value = 0.00155006418458712
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.36964133333333
$ws.Range("N4").Value = 76.108924
$ws.Range("O4").Value = 0.2337290805561598
$ws.Range("P4").Value = 0.2337290805561598
$ws.Range("Q4").Value = 231.0684606823462
$ws.Range("R4").Value = 2079.616146141116
$ws.Range("S4").Value = 0.000362295076666581
$ws.Range("T4").Value = 0.000362295076666581

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.108069666666667
$ws.Range("H5").Value = 27.324209
$ws.Range("I5").Value = 0.00155006418458712
$ws.Range("J5").Value = 0.00155006418458712
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.49855033333333
$ws.Range("N5").Value = 34.495651
$ws.Range("O5").Value = 0.1059354983315251
$ws.Range("P5").Value = 0.1059354983315251
$ws.Range("Q5").Value = 104.7295975016732
$ws.Range("R5").Value = 942.566377515059
$ws.Range("S5").Value = 0.0001642068218400857
$ws.Range("T5").Value = 0.0001642068218400857

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5771.873535333333
$ws.Range("H6").Value = 17315.620606
$ws.Range("I6").Value = 0.9822909543423312
$ws.Range("J6").Value = 0.9822909543423313
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.22896466666667
$ws.Range("N6").Value = 93.686894
$ws.Range("O6").Value = 0.2877106972998646
$ws.Range("P6").Value = 0.2877106972998646
$ws.Range("Q6").Value = 180249.6346953931
$ws.Range("R6").Value = 1622246.712258538
$ws.Range("S6").Value = 0.2826156154251815
$ws.Range("T6").Value = 0.2826156154251816

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5771.873535333333
$ws.Range("H7").Value = 17315.620606
$ws.Range("I7").Value = 0.9822909543423312
$ws.Range("J7").Value = 0.9822909543423313
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 40.44578266666667
$ws.Range("N7").Value = 121.337348
$ws.Range("O7").Value = 0.3726247238124506
$ws.Range("P7").Value = 0.3726247238124505
$ws.Range("Q7").Value = 233447.942589577
$ws.Range("R7").Value = 2101031.483306193
$ws.Range("S7").Value = 0.3660258955652796
$ws.Range("T7").Value = 0.3660258955652796

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5771.873535333333
$ws.Range("H8").Value = 17315.620606
$ws.Range("I8").Value = 0.9822909543423312
$ws.Range("J8").Value = 0.9822909543423313
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.36964133333333
$ws.Range("N8").Value = 76.108924
$ws.Range("O8").Value = 0.2337290805561598
$ws.Range("P8").Value = 0.2337290805561598
$ws.Range("Q8").Value = 146430.3614127653
$ws.Range("R8").Value = 1317873.252714888
$ws.Range("S8").Value = 0.2295899615970658
$ws.Range("T8").Value = 0.2295899615970658

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5771.873535333333
$ws.Range("H9").Value = 17315.620606
$ws.Range("I9").Value = 0.9822909543423312
$ws.Range("J9").Value = 0.9822909543423313
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.49855033333333
$ws.Range("N9").Value = 34.495651
$ws.Range("O9").Value = 0.1059354983315251
$ws.Range("P9").Value = 0.1059354983315251
$ws.Range("Q9").Value = 66368.17836366495
$ws.Range("R9").Value = 597313.6052729846
$ws.Range("S9").Value = 0.1040594817548043
$ws.Range("T9").Value = 0.1040594817548043

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.272029666666667
$ws.Range("H10").Value = 3.816089
$ws.Range("I10").Value = 0.0002164813950916887
$ws.Range("J10").Value = 0.0002164813950916887
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.22896466666667
$ws.Range("N10").Value = 93.686894
$ws.Range("O10").Value = 0.2877106972998646
$ws.Range("P10").Value = 0.2877106972998646
$ws.Range("Q10").Value = 39.72416951528511
$ws.Range("R10").Value = 357.5175256375659
$ws.Range("S10").Value = 0.00006228401313427724
$ws.Range("T10").Value = 0.00006228401313427725

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.272029666666667
$ws.Range("H11").Value = 3.816089
$ws.Range("I11").Value = 0.0002164813950916887
$ws.Range("J11").Value = 0.0002164813950916887
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 40.44578266666667
$ws.Range("N11").Value = 121.337348
$ws.Range("O11").Value = 0.3726247238124506
$ws.Range("P11").Value = 0.3726247238124505
$ws.Range("Q11").Value = 51.44823544355244
$ws.Range("R11").Value = 463.0341189919719
$ws.Range("S11").Value = 0.00008066632005657449
$ws.Range("T11").Value = 0.00008066632005657449

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.272029666666667
$ws.Range("H12").Value = 3.816089
$ws.Range("I12").Value = 0.0002164813950916887
$ws.Range("J12").Value = 0.0002164813950916887
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 25.36964133333333
$ws.Range("N12").Value = 76.108924
$ws.Range("O12").Value = 0.2337290805561598
$ws.Range("P12").Value = 0.2337290805561598
$ws.Range("Q12").Value = 32.27093640869289
$ws.Range("R12").Value = 290.438427678236
$ws.Range("S12").Value = 0.00005059799743229515
$ws.Range("T12").Value = 0.00005059799743229516

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.272029666666667
$ws.Range("H13").Value = 3.816089
$ws.Range("I13").Value = 0.0002164813950916887
$ws.Range("J13").Value = 0.0002164813950916887
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 11.49855033333333
$ws.Range("N13").Value = 34.495651
$ws.Range("O13").Value = 0.1059354983315251
$ws.Range("P13").Value = 0.1059354983315251
$ws.Range("Q13").Value = 14.62649714765989
$ws.Range("R13").Value = 131.638474328939
$ws.Range("S13").Value = 0.00002293306446854183
$ws.Range("T13").Value = 0.00002293306446854183

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 93.67702500000001
$ws.Range("H14").Value = 281.031075
$ws.Range("I14").Value = 0.01594250007799006
$ws.Range("J14").Value = 0.01594250007799006
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.22896466666667
$ws.Range("N14").Value = 93.686894
$ws.Range("O14").Value = 0.2877106972998646
$ws.Range("P14").Value = 0.2877106972998646
$ws.Range("Q14").Value = 2925.43650380345
$ws.Range("R14").Value = 26328.92853423105
$ws.Range("S14").Value = 0.004586827814141666
$ws.Range("T14").Value = 0.004586827814141666

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 93.67702500000001
$ws.Range("H15").Value = 281.031075
$ws.Range("I15").Value = 0.01594250007799006
$ws.Range("J15").Value = 0.01594250007799006
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 40.44578266666667
$ws.Range("N15").Value = 121.337348
$ws.Range("O15").Value = 0.3726247238124506
$ws.Range("P15").Value = 0.3726247238124505
$ws.Range("Q15").Value = 3788.840594009901
$ws.Range("R15").Value = 34099.5653460891
$ws.Range("S15").Value = 0.005940569688441018
$ws.Range("T15").Value = 0.005940569688441017

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 93.67702500000001
$ws.Range("H16").Value = 281.031075
$ws.Range("I16").Value = 0.01594250007799006
$ws.Range("J16").Value = 0.01594250007799006
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 25.36964133333333
$ws.Range("N16").Value = 76.108924
$ws.Range("O16").Value = 0.2337290805561598
$ws.Range("P16").Value = 0.2337290805561598
$ws.Range("Q16").Value = 2376.5525254237
$ws.Range("R16").Value = 21388.9727288133
$ws.Range("S16").Value = 0.003726225884995122
$ws.Range("T16").Value = 0.003726225884995122

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 93.67702500000001
$ws.Range("H17").Value = 281.031075
$ws.Range("I17").Value = 0.01594250007799006
$ws.Range("J17").Value = 0.01594250007799006
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 11.49855033333333
$ws.Range("N17").Value = 34.495651
$ws.Range("O17").Value = 0.1059354983315251
$ws.Range("P17").Value = 0.1059354983315251
$ws.Range("Q17").Value = 1077.149987039425
$ws.Range("R17").Value = 9694.349883354827
$ws.Range("S17").Value = 0.001688876690412256
$ws.Range("T17").Value = 0.001688876690412256

